# Applies the "LinuxForHealth" re-brand edit to the
# StructureDefinition-match-confidence-score workbook:
#   - Metadata sheet: URL / Version / Date / Publisher values updated
#   - Elements sheet: the erroneous duplicate ele-1/ext-1 constraint text
#     that was (wrongly) shown on the "Extension" row is cleared, leaving
#     it only on the "Extension.extension" row where it belongs. The
#     "Fixed Value" for Extension.url is refreshed to the new base URL.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-score"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Clear the stray ele-1/ext-1 constraint text from the top-level
# "Extension" row (row 2) - it belongs only to "Extension.extension" (row 4).
$elements.Range("AI2").Value = ""

# Keep the "Fixed Value" for Extension.url in sync with the new base URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-score"
